# electricalDiagram.pptx — "adding programming and software files"
#
# 1) The deck's cached datetimeFigureOut field ("today" stamped into the
#    Date placeholder on the slide master and every slide layout) moves
#    forward one day: 2/22/2016 -> 2/23/2016.
# 2) Three shapes on slide 1 shift slightly to the right (X offset only;
#    Y/width/height are untouched):
#      - "Straight Connector 96"  (rotated connector)
#      - "Straight Connector 170"
#      - "Oval 172"

$p = $ppt.ActivePresentation

# EMU -> point helper (PowerPoint's Shape.Left/.Top are in points).
function EmuToPt([double]$emu) { return $emu / 12700.0 }

# --- 1) Refresh the cached date-field text everywhere it appears -------

$oldDate = "2/22/2016"
$newDate = "2/23/2016"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout hanging off the master.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# --- 2) Nudge the three connector/oval shapes on slide 1 ---------------

$slide = $p.Slides.Item(1)
$shapes = $slide.Shapes

$moves = @{
    "Straight Connector 96"  = EmuToPt 3360901
    "Straight Connector 170" = EmuToPt 3382497
    "Oval 172"                = EmuToPt 3519551
}

for ($i = 1; $i -le $shapes.Count; $i++) {
    $sh = $shapes.Item($i)
    if ($moves.ContainsKey($sh.Name)) {
        $sh.Left = $moves[$sh.Name]
    }
}
